$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '46.743.12'
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  +0.10%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.272.85'
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  -3.00%  '
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '299.76'
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  -2.38%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '99.41'
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  +0.12%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.573'
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  -0.92%  '
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.507'
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  -5.44%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '34.99'
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  -2.83%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0802'
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '7.02'
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  -5.46%  '
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  -1.74%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '2.622.31'
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  -3.03%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '2.275.36'
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  -2.96%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '13.62'
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  -4.60%  '
$ws.Range('B17').NumberFormat = "@"
$ws.Range('B17').Value = 'WrappedBTC'
$ws.Range('C17').NumberFormat = "@"
$ws.Range('C17').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '46.686.21'
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  -0.10%  '
$ws.Range('B18').NumberFormat = "@"
$ws.Range('B18').Value = 'Polygon'
$ws.Range('C18').NumberFormat = "@"
$ws.Range('C18').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.798'
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  -3.82%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.0₃0990'
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  +4.35%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '12.48'
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  -6.66%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '5.81'
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  -5.90%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '65.84'
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  -1.30%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '246.82'
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  +0.63%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.78'
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  -6.52%  '
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  +0.18%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '1.85'
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  -6.48%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '40.98'
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  -2.29%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '2.21'
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  -3.37%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '9.54'
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  -3.49%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '20.07'
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  -0.76%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '2.81'
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  +7.43%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.34'
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  +8.63%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '146.26'
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  -3.37%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '5.32'
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  -7.11%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.0762'
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  -6.20%  '
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  +3.31%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.115'
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  -2.70%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '15.52'
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  +12.03%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '1.67'
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  -8.87%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '3.83'
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  -5.04%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.0296'
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  -7.35%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '3.07'
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  -10.65%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.999'
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  -0.16%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '93.74'
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  +15.56%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '1.784.14'
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  -1.27%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '1.87'
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  -3.56%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '70.53'
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  -3.67%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.184'
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  -7.05%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '4.80'
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  -2.65%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '94.49'
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  -4.03%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '7.88'
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  -1.48%  '
